{"js": "// The document talks about the \"pom.xml.jim\" Jamal template file; the\n// author renamed/retyped it in the sentence \"This is the content of the\n// pom.xml.jim file:\" so that it now reads \"pom.jam\" (the template's new\n// name), keeping the existing Courier New / 9pt / en-US run formatting\n// that was already applied to that piece of text.\nconst results = context.document.body.search(\"pom.xml.jim\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'pom.xml.jim' in the document body.\");\n}\n\n// Replace the found text in place; Word keeps the run's existing\n// character formatting (font, size, language) for the replacement text.\nresults.items[0].insertText(\"pom.jam\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document talks about the \"pom.xml.jim\" Jamal template file; the\n# author renamed/retyped it in the sentence \"This is the content of the\n# pom.xml.jim file:\" so that it now reads \"pom.jam\" (the template's new\n# name). Use Find/Replace on the document body so Word preserves the\n# existing run formatting (Courier New, 9pt, en-US) already applied to\n# that piece of text.\n$d = $word.ActiveDocument\n$rng = $d.Content\n\n$rng.Find.ClearFormatting()\n$rng.Find.Replacement.ClearFormatting()\n\n$found = $rng.Find.Execute(\"pom.xml.jim\", $true, $false, $false, $false, $false, $true, 1, $false, \"pom.jam\", 2)\n\nif (-not $found) {\n    throw \"Could not find 'pom.xml.jim' in the document.\"\n}\n"}
